# BOM.xlsx update: add Price / Total columns, Final Total and
# "What we have to pay" summary rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the old "datasheet link" column (F) one column to the right
# (to G) so the new Price/Total columns (D/E) can be inserted before it.
$ws.Columns("F:F").Insert()

# --- Column widths for the new D/E (Price/Total) and the (now) G hyperlink column
$ws.Columns("D:E").ColumnWidth = 19.85546875
$ws.Columns("F:F").ColumnWidth = 18.5703125

# --- Header row
$ws.Range("D1").Value = "Price"
$ws.Range("E1").Value = "Total"

# --- Per-part unit price (D) and extended total (E = Qty * Price)
$prices = @{
    2  = 14.84
    3  = 6.8
    4  = 6.8
    5  = 0.48
    6  = 0.68
    7  = 1.13
    8  = 2.58
    9  = 1.87
    10 = 0.25
    11 = 3.35
    12 = 6.26
    13 = 0.69
    14 = 0.74
    15 = 6.42
    16 = 2.37
    17 = 0.41
    18 = 3.48
    19 = 3.39
    20 = 0.86
    21 = 8.67
    22 = 0.21
    23 = 2.68
    24 = 1.2
    25 = 0.88
    26 = 0.12
    27 = 0.59
    28 = 0.61
    29 = 0.29
    30 = 0.29
}

foreach ($row in 2..30) {
    $ws.Cells.Item($row, 4).Value = $prices[$row]
}

$ws.Range("E2").Formula = "=`$C2*`$D2"
$ws.Range("E3:E30").Formula = "=`$C3*`$D3"

$ws.Range("D2:E30").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

# --- Summary rows
$ws.Range("E32").Value = "Final Total:"
$ws.Range("F32").Formula = "=SUM(`$E2:`$E30)"

$ws.Range("E33").Value = "What we have to pay:"
$ws.Range("F33").Formula = "=`$F`$32-`$E`$11-`$E`$12-`$E`$23-`$E`$24"

$ws.Range("F32:F33").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

# --- Selection the workbook ended on
$ws.Range("G35").Select()
